$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# The title currently holds three runs ("Lists", " ", "(continued)") whose
# concatenation already equals the desired text, so a direct assignment of
# the same combined string is treated as a same-value write and leaves the
# runs untouched. Force the consolidation by nudging the text through a
# temporary value first, then set the final desired text which rewrites the
# paragraph as a single run.
$tr.Text = "__tmp__"
$tr.Text = "Lists (continued)"
